$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.291.67'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.494.55'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.69'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.72'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +6.34%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.387'
$ws.Range('E11').Value = '  +2.75%  '
$ws.Range('D12').Value = '4.089.55'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '3.493.87'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '64.224.62'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.31'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.06'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('D23').Value = '3.633.51'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.11'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('E33').Value = '  +3.79%  '
$ws.Range('D34').Value = '3.523.64'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.58'
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0784'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.807'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.45'
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D48').Value = '2.436.17'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.81'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.913'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('E51').Value = '  -0.62%  '
